$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("30-Dec-2023")

# Re-order/update the first four data rows (rows 2-5) to reflect the
# refreshed sort order + updated NTT/LTT values for the ARE&M.NS stock.
$ws.Range("A2").Value = "ARE&M.NS"
$ws.Range("B2").Value = 20.47
$ws.Range("C2").Value = 20.47
$ws.Range("D2").Value = 2.77
$ws.Range("E2").Value = "NTT"

$ws.Range("A3").Value = "SAIL.NS"
$ws.Range("B3").Value = 4.9400000000000004
$ws.Range("C3").Value = 4.9400000000000004
$ws.Range("D3").Value = 0.71
$ws.Range("E3").Value = "ATH"

$ws.Range("A4").Value = "EMAMILTD.NS"
$ws.Range("B4").Value = 5.3
$ws.Range("C4").Value = 5.3
$ws.Range("D4").Value = 0.65
$ws.Range("E4").Value = "ATH"

$ws.Range("A5").Value = "PIDILITIND.NS"
$ws.Range("B5").Value = 6.43
$ws.Range("C5").Value = 6.74
$ws.Range("D5").Value = 1.67
$ws.Range("E5").Value = "BTT"

# Updated NTT% figures further down the table (refreshed stock data).
$ws.Range("B27").Value = 8.6199999999999992
$ws.Range("B44").Value = 13.09
$ws.Range("B51").Value = 17.88
$ws.Range("B56").Value = 19.05
$ws.Range("B57").Value = 4.68

# Update the last active selection to match the author's saved state.
$ws.Range("I15").Select()
